$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column price values are written as text (matches original inlineStr type),
# since some values (e.g. "1.00", "0.514") would otherwise be auto-coerced to numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "41.524.08"
$ws.Range("D3").Value = "2.468.23"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "314.81"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").Value = "91.93"
$ws.Range("E6").Value = "  -2.69%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("D9").Value = "0.514"
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("D10").Value = "32.30"
$ws.Range("E10").Value = "  -3.92%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "2.847.64"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "15.95"
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("D16").Value = "2.459.19"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").Value = "0.775"
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").Value = "41.522.72"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "6.47"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("D21").Value = "71.07"
$ws.Range("E21").Value = "  +3.18%  "
$ws.Range("D22").Value = "11.07"
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").Value = "236.32"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "2.71"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  -0.69%  "
$ws.Range("D27").Value = "24.65"
$ws.Range("E27").Value = "  +1.90%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").Value = "35.37"
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("D31").Value = "155.67"
$ws.Range("E31").Value = "  +2.18%  "
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "17.23"
$ws.Range("E35").Value = "  -4.84%  "
$ws.Range("E36").Value = "  -7.02%  "
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").Value = "  -5.71%  "
$ws.Range("E40").Value = "  -11.00%  "
$ws.Range("E41").Value = "  -5.05%  "
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").Value = "1.945.75"
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").Value = "18.59"
$ws.Range("E45").Value = "  -6.26%  "
$ws.Range("E46").Value = "  -3.31%  "
$ws.Range("D47").Value = "9.04"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("D48").Value = "2.706.04"
$ws.Range("E48").Value = "  -1.06%  "
$ws.Range("D49").Value = "96.90"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").Value = "67.18"
$ws.Range("E50").Value = "  -3.95%  "
$ws.Range("D51").Value = "52.04"
$ws.Range("E51").Value = "  +1.92%  "

# Restore default style (removes the temporary text-number-format override)
$priceRange.Style = "Normal"
